$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sponza" (sheet2.xml): add two new result columns N (v1436) / O (v1437)
# ---------------------------------------------------------------------------
$sponza = $wb.Worksheets.Item("Sponza")

# Carry formatting (header/data/average/variance styles) from column M into N
# and O before writing any values/formulas, so every new cell inherits the
# same style index as its neighbour in the same row.
$sponza.Range("M1:M16").Copy()
$sponza.Range("N1:N16").PasteSpecial(-4122)
$sponza.Range("N1:N16").Copy()
$sponza.Range("O1:O16").PasteSpecial(-4122)

# Header labels (new shared strings "v1436" / "v1437")
$sponza.Cells.Item(1, 14).Value = "v1436"
$sponza.Cells.Item(1, 15).Value = "v1437"

# Raw per-run samples
$sponzaN = @(6260,6243,6262,6259,6246,6257,6277,6237,6242,6274)
$sponzaO = @(5809,5794,5789,5790,5802,5796,5765,5767,5792,5783)
for ($i = 0; $i -lt 10; $i++) {
    $sponza.Cells.Item(2 + $i, 14).Value = $sponzaN[$i]
    $sponza.Cells.Item(2 + $i, 15).Value = $sponzaO[$i]
}

# AVG / VAR / DIFF ACCEPT / ratio rows, same formula shapes as column M
$sponza.Range("N12:O12").Formula = "=AVERAGE(N2:N11)"
$sponza.Range("N13:O13").Formula = "=_xlfn.VAR.S(N2:N11)"
$sponza.Range("N14:O14").Formula = "=1-_xlfn.T.TEST(M2:M11,N2:N11,2,3)"
$sponza.Range("N15:O15").Formula = "=M12/N12"
$sponza.Cells.Item(16, 14).Formula = "=B12/N12"
$sponza.Cells.Item(16, 15).Formula = "=B12/O12"

# Conditional formatting used to cover B15:M16; extend it to cover the two
# new columns as well.
$fmtConds = $sponza.Range("B15:M16").FormatConditions
for ($i = 1; $i -le $fmtConds.Count; $i++) {
    $fmtConds.Item($i).ModifyAppliesToRange($sponza.Range("B15:O16"))
}

# ---------------------------------------------------------------------------
# Sheet "ComplexMesh" (sheet3.xml): add two new result columns M (v1436) /
# N (v1437); also turn the previously-hardcoded H12:L12 averages into live
# AVERAGE() formulas (matching the rest of the row).
# ---------------------------------------------------------------------------
$complexMesh = $wb.Worksheets.Item("ComplexMesh")

$complexMesh.Range("C12").Copy()
$complexMesh.Range("H12:L12").PasteSpecial(-4122)
$complexMesh.Range("H12:L12").Formula = "=AVERAGE(H2:H11)"

$complexMesh.Range("L1:L16").Copy()
$complexMesh.Range("M1:M16").PasteSpecial(-4122)
$complexMesh.Range("M1:M16").Copy()
$complexMesh.Range("N1:N16").PasteSpecial(-4122)

$complexMesh.Cells.Item(1, 13).Value = "v1436"
$complexMesh.Cells.Item(1, 14).Value = "v1437"

$complexMeshM = @(4187,4221,4219,4238,4168,4183,4180,4157,4144,4167)
$complexMeshN = @(3968,3961,3964,3950,3980,3944,3953,3933,3947,3959)
for ($i = 0; $i -lt 10; $i++) {
    $complexMesh.Cells.Item(2 + $i, 13).Value = $complexMeshM[$i]
    $complexMesh.Cells.Item(2 + $i, 14).Value = $complexMeshN[$i]
}

$complexMesh.Range("M12:N12").Formula = "=AVERAGE(M2:M11)"
$complexMesh.Range("M13:N13").Formula = "=_xlfn.VAR.S(M2:M11)"
$complexMesh.Range("M14:N14").Formula = "=1-_xlfn.T.TEST(L2:L11,M2:M11,2,3)"
$complexMesh.Range("M15:N15").Formula = "=L12/M12"
$complexMesh.Cells.Item(16, 13).Formula = "=B12/M12"
$complexMesh.Cells.Item(16, 14).Formula = "=B12/N12"

$fmtConds3 = $complexMesh.Range("B15:L16").FormatConditions
for ($i = 1; $i -le $fmtConds3.Count; $i++) {
    $fmtConds3.Item($i).ModifyAppliesToRange($complexMesh.Range("B15:N16"))
}

# ---------------------------------------------------------------------------
# Selections per-sheet + which tab is active.
# ---------------------------------------------------------------------------
$partOfSponza = $wb.Worksheets.Item("PartOfSponza")
$partOfSponza.Activate()
$partOfSponza.Range("G1").Select()

$complexMesh.Activate()
$complexMesh.Range("M7").Select()

$sponza.Activate()
$sponza.Range("O2").Select()
